$wb = $excel.ActiveWorkbook

# Work on the "Backlog" sheet: mark US04 (Marriage before divorce) and
# US05 (Marriage before death) statuses as Complete.
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("E4").Value = "Complete"
$backlog.Range("E5").Value = "Complete"
$backlog.Activate()
[void]$backlog.Range("F4").Select()

# Browse through the other sprint-tracking sheets (selection/navigation
# only, no data changes on these).
$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Activate()
[void]$burndown.Range("C32").Select()

$stories = $wb.Worksheets.Item("Stories")
$stories.Activate()
[void]$stories.Range("C6").Select()

# Work on the "Sprint1" sheet: mark US04 (Marriage before divorce) and
# US05 (Marriage before death) as complete.
$ws = $wb.Worksheets.Item("Sprint1")
$ws.Activate()

# Row 4 -> US04 "Marriage before divorce" (owner MW)
$ws.Range("D4").Value = "Complete"
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = "Yes"

# Row 5 -> US05 "Marriage before death" (owner MW)
$ws.Range("D5").Value = "Complete"
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = "Yes"

[void]$ws.Range("B5").Select()
